$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 114 (Clemenuless/Especial),
# pushing the existing rows 114-162 down to 116-164.
$ws.Rows.Item(114).Resize(2).Insert()

# New row 114: Murcott / Primera entry for 2022-10-11 (serial 44845)
$ws.Range("A114").Value = 11
$ws.Range("B114").Value = "Vega Monumental Concepción"
$ws.Range("C114").Value = "Bíobío"
$ws.Range("D114").Value = 44845
$ws.Range("E114").Value = 8
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100102
$ws.Range("H114").Value = "Cítricos"
$ws.Range("I114").Value = 100102004
$ws.Range("J114").Value = "Mandarina"
$ws.Range("K114").Value = "Murcott"
$ws.Range("L114").Value = "Primera"
$ws.Range("M114").Value = 200
$ws.Range("N114").Value = 6000
$ws.Range("O114").Value = 6500
$ws.Range("P114").Value = 6250
$ws.Range("Q114").Value = "`$/bandeja 18 kilos"
$ws.Range("R114").Value = "Región de O'Higgins"
$ws.Range("S114").Value = 347
$ws.Range("T114").Value = 18

# New row 115: Murcott / Segunda entry for the same date
$ws.Range("A115").Value = 11
$ws.Range("B115").Value = "Vega Monumental Concepción"
$ws.Range("C115").Value = "Bíobío"
$ws.Range("D115").Value = 44845
$ws.Range("E115").Value = 8
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100102
$ws.Range("H115").Value = "Cítricos"
$ws.Range("I115").Value = 100102004
$ws.Range("J115").Value = "Mandarina"
$ws.Range("K115").Value = "Murcott"
$ws.Range("L115").Value = "Segunda"
$ws.Range("M115").Value = 100
$ws.Range("N115").Value = 5000
$ws.Range("O115").Value = 5000
$ws.Range("P115").Value = 5000
$ws.Range("Q115").Value = "`$/bandeja 18 kilos"
$ws.Range("R115").Value = "Región de O'Higgins"
$ws.Range("S115").Value = 278
$ws.Range("T115").Value = 18
